$wb = $excel.ActiveWorkbook

# --- DatosCuenta: update account holder name/lastname + documento + numero calle ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokeCatorce"
$wsCuenta.Range("B2").Value = "SmokeLastNCatorce"
$wsCuenta.Range("C2").Value = 20111101
$wsCuenta.Range("D2").Value = 101

# --- DatosHogar: bump household number ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 621

# --- DatosMotor: update vehicle plate/motor/chassis numbers ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMA002"
$wsMotor.Range("B2").Value = "ABC12SSMA002"
$wsMotor.Range("C2").Value = "ZAZ123SSMA002"

# --- DatosAP: bump AP number ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200101

# Select A2:C2 on DatosMotor (matches the saved selection in the workbook)
# without leaving that sheet as the active tab -- DatosAP stays the active tab.
$wsMotor.Range("A2:C2").Select()
$wsAP.Activate()
